$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (17-nov) before DV ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at DV; existing DV:EZ content shifts right to DW:FA
$ws1.Range("DV1").EntireColumn.Insert()

# Header for the newly inserted column
$ws1.Range("DV1").Value = "17-nov"

# Data rows 2-25 get a placeholder "-" in the new column, same as the other
# not-yet-available days
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 126).Value = "-"
}

# --- Sheet "Gaz": append two more daily price rows ---
$ws2 = $wb.Worksheets.Item("Gaz")

$ws2.Range("A153").NumberFormat = "@"
$ws2.Range("A153").Value = "2025-11-15"
$ws2.Range("A153").Style = "Normal"
$ws2.Range("B153").Value = 29.305

$ws2.Range("A154").NumberFormat = "@"
$ws2.Range("A154").Value = "2025-11-16"
$ws2.Range("A154").Style = "Normal"
$ws2.Range("B154").Value = 29.305

# --- Sheet "CO2": append two more daily price rows ---
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A154").NumberFormat = "@"
$ws3.Range("A154").Value = "2025-11-15"
$ws3.Range("A154").Style = "Normal"
$ws3.Range("B154").Value = 80.72

$ws3.Range("A155").NumberFormat = "@"
$ws3.Range("A155").Value = "2025-11-16"
$ws3.Range("A155").Style = "Normal"
$ws3.Range("B155").Value = 80.72
